$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction7")

# Update A1 and B1 with the new random, non-overlapping values
$ws.Range("A1").Value = 14
$ws.Range("B1").Value = 15

# Remove the remaining values C1:O1 so the sheet only spans A1:B1
$ws.Range("C1:O1").Clear()
